# Refresh the cached "today" text shown by the auto-updating Date
# placeholders on the slide master, slide layouts, notes master and
# handout master (Insert > Header & Footer > Date and time > Update
# automatically). PowerPoint recomputes these caches whenever the
# deck is saved; this mirrors that re-cache after the file was
# reopened/saved on 12/15/2024 (most masters/layouts + notes master)
# and 12/16/2024 (handout master).

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide master date placeholder -> 12/15/2024
Set-DatePlaceholderText $p.SlideMaster.Shapes "12/15/2024"

# Slide layouts date placeholders -> 12/15/2024 (layouts without a
# date placeholder are silently skipped by Set-DatePlaceholderText)
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes "12/15/2024"
}

# Notes master date placeholder -> 12/15/2024
Set-DatePlaceholderText $p.NotesMaster.Shapes "12/15/2024"

# Handout master date placeholder -> 12/16/2024
Set-DatePlaceholderText $p.HandoutMaster.Shapes "12/16/2024"
